$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.918.10'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '1.555.38'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.486'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.39%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.247'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.51'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0859'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("D12").Value = '1.777.47'
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Value = '1.557.83'
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("E14").Value = '  +0.64%  '
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").Value = '26.914.28'
$ws.Range("E16").Value = '  +0.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.66'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").Value = '0.0₃0687'
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("E22").Value = '  -0.92%  '
$ws.Range("E23").Value = '  +1.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("E26").Value = '  +2.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("E32").Value = '  +2.21%  '
$ws.Range("D33").Value = '1.368.62'
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("E34").Value = '  +1.80%  '
$ws.Range("E35").Value = '  +3.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.967'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.99%  '
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("E38").Value = '  +0.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.522'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.808'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.87%  '
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.992'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("E43").Value = '  -0.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.69'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.05%  '
$ws.Range("E46").Value = '  -1.38%  '
$ws.Range("D47").Value = '1.690.41'
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("E49").Value = '  -1.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0954'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.34%  '
